$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 5247000000.0
$ws.Range("C4").Value = 5026000000.0
$ws.Range("D4").Value = 4391000000.0
$ws.Range("E4").Value = 4107000000.0
$ws.Range("F4").Value = 4677000000.0

# Row 12 - Accounts Payable
$ws.Range("B12").Value = 3614000000.0
$ws.Range("C12").Value = 3771000000.0
$ws.Range("D12").Value = 3401000000.0
$ws.Range("E12").Value = 2954000000.0
$ws.Range("F12").Value = 2861000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("B20").Value = 711000000.0
$ws.Range("C20").Value = 687000000.0
$ws.Range("D20").Value = 690000000.0
$ws.Range("E20").Value = 700000000.0
$ws.Range("F20").Value = 675000000.0
